$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove now-unused rows 56-84 (data was deduplicated/cleaned)
$ws.Range("A56:C84").EntireRow.Delete()

# Rewrite A1:C55 with the deduplicated, relabeled, sorted data
$ws.Range("A1").Value2 = "hebrew_text"
$ws.Range("B1").Value2 = "hebrew_option"
$ws.Range("C1").Value2 = "label"
$ws.Range("A2").Value2 = "מאז הסקר האחרון, דעתי הוסחה בקלות"
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value2 = "C_ADHD_Distracted"
$ws.Range("A3").Value2 = "ברגע זה אני מרגיש חסר מנוחה"
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value2 = "C_ADHD_Restless"
$ws.Range("A4").Value2 = "ברגע זה, אני מרגיש חסר מנוחה"
$ws.Range("B4").ClearContents()
$ws.Range("C4").Value2 = "C_ADHD_Restless"
$ws.Range("A5").Value2 = "ברגע זה, אני מרגישה חסרת מנוחה"
$ws.Range("B5").ClearContents()
$ws.Range("C5").Value2 = "C_ADHD_Restless"
$ws.Range("A6").Value2 = "ברגע זה, אני מרגישה חסר מנוחה"
$ws.Range("B6").ClearContents()
$ws.Range("C6").Value2 = "C_ADHD_Restless"
$ws.Range("A7").Value2 = "מאז הסקר האחרון, התפרצתי בכעס באחת או יותר מהדרכים הבאות"
$ws.Range("B7").Value2 = "הרבצתי"
$ws.Range("C7").Value2 = "C_Agr_hit"
$ws.Range("A8").Value2 = "מאז הסקר האחרון, התפרצתי בכעס באחת או יותר מהדרכים הבאות"
$ws.Range("B8").Value2 = "לא התפרצתי בכעס"
$ws.Range("C8").Value2 = "C_Agr_none"
$ws.Range("A9").Value2 = "מאז הסקר האחרון, כעסתי או התעצבנתי כשדברים קרו לא כמו שרציתי"
$ws.Range("B9").ClearContents()
$ws.Range("C9").Value2 = "C_Agr_NotAsWant"
$ws.Range("A10").Value2 = "מאז הסקר האחרון, התפרצתי בכעס באחת או יותר מהדרכים הבאות"
$ws.Range("B10").Value2 = "התפרצתי בכעס בצורה אחרת (בבקשה לפרט)"
$ws.Range("C10").Value2 = "C_Agr_other"
$ws.Range("A11").Value2 = "מאז הסקר האחרון, התפרצתי בכעס באחת או יותר מהדרכים הבאות"
$ws.Range("B11").Value2 = "טרקתי דלת"
$ws.Range("C11").Value2 = "C_Agr_slam"
$ws.Range("A12").Value2 = "מאז הסקר האחרון, התפרצתי בכעס באחת או יותר מהדרכים הבאות"
$ws.Range("B12").Value2 = "זרקתי משהו"
$ws.Range("C12").Value2 = "C_Agr_throw_smt"
$ws.Range("A13").Value2 = "מאז הסקר האחרון, התפרצתי בכעס באחת או יותר מהדרכים הבאות"
$ws.Range("B13").Value2 = "זרקתי משהו על מישהו"
$ws.Range("C13").Value2 = "C_Agr_throw_twd"
$ws.Range("A14").Value2 = "מאז הסקר האחרון, התפרצתי בכעס באחת או יותר מהדרכים הבאות"
$ws.Range("B14").Value2 = "צעקתי"
$ws.Range("C14").Value2 = "C_Agr_yelled"
$ws.Range("A15").Value2 = "ברגע זה, אני מרגיש מרוגז או כועס"
$ws.Range("B15").ClearContents()
$ws.Range("C15").Value2 = "C_Angry_now"
$ws.Range("A16").Value2 = "ברגע זה, אני מרגישה מרוגזת או כועסת"
$ws.Range("B16").ClearContents()
$ws.Range("C16").Value2 = "C_Angry_now"
$ws.Range("A17").Value2 = "ברגע זה, אני מרגיש פחד או לחץ"
$ws.Range("B17").ClearContents()
$ws.Range("C17").Value2 = "C_Anx_now"
$ws.Range("A18").Value2 = "ברגע זה, אני מרגישה פחד או לחץ"
$ws.Range("B18").ClearContents()
$ws.Range("C18").Value2 = "C_Anx_now"
$ws.Range("A19").Value2 = "מאז הסקר האחרון, הרגשתי מודאג או מפוחד"
$ws.Range("B19").ClearContents()
$ws.Range("C19").Value2 = "C_Anx_Worry"
$ws.Range("A20").Value2 = "מאז הסקר האחרון, הרגשתי מודאגת או מפוחדת"
$ws.Range("B20").ClearContents()
$ws.Range("C20").Value2 = "C_Anx_Worry"
$ws.Range("A21").Value2 = "מאז הסקר האחרון, אבא שלי איים להעניש אותי אבל לא עשה זאת"
$ws.Range("B21").ClearContents()
$ws.Range("C21").Value2 = "C_Discipline"
$ws.Range("A22").Value2 = "מאז הסקר האחרון, אמא שלי איימה להעניש אותי אבל לא עשתה זאת"
$ws.Range("B22").ClearContents()
$ws.Range("C22").Value2 = "C_Discipline"
$ws.Range("A23").Value2 = "מאז הסקר האחרון, היה לי קשה להפסיק לעשות משהו אחרי שביקשו ממני להפסיק"
$ws.Range("B23").ClearContents()
$ws.Range("C23").Value2 = "C_IC_CantStop"
$ws.Range("A24").Value2 = "מאז הסקר האחרון, יצא שאמרתי את הדבר הראשון שעלה לי לראש מבלי לעצור ולחשוב"
$ws.Range("B24").ClearContents()
$ws.Range("C24").Value2 = "C_IC_FirstOnMind"
$ws.Range("A25").Value2 = "היום אני ואבא שלי בילינו זמן כיף ביחד"
$ws.Range("B25").ClearContents()
$ws.Range("C25").Value2 = "C_Inv_Fun"
$ws.Range("A26").Value2 = "היום אני ואמא שלי בילינו זמן כיף ביחד"
$ws.Range("B26").ClearContents()
$ws.Range("C26").Value2 = "C_Inv_Fun"
$ws.Range("A27").Value2 = "היום אבא שלי עזר לי במשהו"
$ws.Range("B27").ClearContents()
$ws.Range("C27").Value2 = "C_Inv_Help"
$ws.Range("A28").Value2 = "היום אמא שלי עזרה לי במשהו"
$ws.Range("B28").ClearContents()
$ws.Range("C28").Value2 = "C_Inv_Help"
$ws.Range("A29").Value2 = "היום אבא שלי דיבר איתי על החברים שלי או על היום שלי"
$ws.Range("B29").ClearContents()
$ws.Range("C29").Value2 = "C_Inv_Talk"
$ws.Range("A30").Value2 = "היום אמא שלי דיברה איתי על החברים שלי או על היום שלי"
$ws.Range("B30").ClearContents()
$ws.Range("C30").Value2 = "C_Inv_Talk"
$ws.Range("A31").Value2 = "מאז הסקר האחרון, הרגשתי מתוסכל"
$ws.Range("B31").ClearContents()
$ws.Range("C31").Value2 = "C_Irr_Frustration"
$ws.Range("A32").Value2 = "מאז הסקר האחרון, הרגשתי מתוסכלת"
$ws.Range("B32").ClearContents()
$ws.Range("C32").Value2 = "C_Irr_Frustration"
$ws.Range("A33").Value2 = "ברגע זה אני מרגיש טוב"
$ws.Range("B33").ClearContents()
$ws.Range("C33").Value2 = "C_Mood_Good"
$ws.Range("A34").Value2 = "ברגע זה, אני מרגיש טוב"
$ws.Range("B34").ClearContents()
$ws.Range("C34").Value2 = "C_Mood_Good"
$ws.Range("A35").Value2 = "ברגע ,זה אני מרגיש טוב"
$ws.Range("B35").ClearContents()
$ws.Range("C35").Value2 = "C_Mood_Good"
$ws.Range("A36").Value2 = "ברגע זה, אני מרגישה טוב"
$ws.Range("B36").ClearContents()
$ws.Range("C36").Value2 = "C_Mood_Good"
$ws.Range("A37").Value2 = "מאז הסקר האחרון, הרגשתי עצוב או מדוכא"
$ws.Range("B37").ClearContents()
$ws.Range("C37").Value2 = "C_Mood_Sad"
$ws.Range("A38").Value2 = "מאז הסקר האחרון, הרגשתי עצובה או מדוכאת"
$ws.Range("B38").ClearContents()
$ws.Range("C38").Value2 = "C_Mood_Sad"
$ws.Range("A39").Value2 = "מאז הסקר האחרון, עצבנתי את אבא שלי"
$ws.Range("B39").ClearContents()
$ws.Range("C39").Value2 = "C_PC_Annoy"
$ws.Range("A40").Value2 = "מאז הסקר האחרון, עצבנתי את אמא שלי"
$ws.Range("B40").ClearContents()
$ws.Range("C40").Value2 = "C_PC_Annoy"
$ws.Range("A41").Value2 = "מאז הסקר האחרון, אבא שלי העביר עליי ביקורת"
$ws.Range("B41").ClearContents()
$ws.Range("C41").Value2 = "C_PC_Criticism"
$ws.Range("A42").Value2 = "מאז הסקר האחרון, אמא שלי העבירה עליי ביקורת"
$ws.Range("B42").ClearContents()
$ws.Range("C42").Value2 = "C_PC_Criticism"
$ws.Range("A43").Value2 = "מאז הסקר האחרון, אבא שלי העביר עליי ביקורת"
$ws.Range("B43").ClearContents()
$ws.Range("C43").Value2 = "C_PC_Criticism"
$ws.Range("A44").Value2 = "מאז הסקר האחרון, שיתפתי את אבא שלי ברגשות/תחושות שלי"
$ws.Range("B44").ClearContents()
$ws.Range("C44").Value2 = "C_PC_Sharing"
$ws.Range("A45").Value2 = "מאז הסקר האחרון, שיתפתי את אמא שלי ברגשות/תחושות שלי"
$ws.Range("B45").ClearContents()
$ws.Range("C45").Value2 = "C_PC_Sharing"
$ws.Range("A46").Value2 = "מאז הסקר האחרון, אבא שלי החמיא לי כשעשיתי משהו בצורה טובה"
$ws.Range("B46").ClearContents()
$ws.Range("C46").Value2 = "C_Positive"
$ws.Range("A47").Value2 = "מאז הסקר האחרון, אמא שלי החמיאה לי כשעשיתי משהו בצורה טובה"
$ws.Range("B47").ClearContents()
$ws.Range("C47").Value2 = "C_Positive"
$ws.Range("A48").Value2 = "מאז הסקר האחרון, אבא שלי הסכים לכל מה שרציתי"
$ws.Range("B48").ClearContents()
$ws.Range("C48").Value2 = "C_PS_Agree"
$ws.Range("A49").Value2 = "מאז הסקר האחרון, אמא שלי הסכימה לכל מה שרציתי"
$ws.Range("B49").ClearContents()
$ws.Range("C49").Value2 = "C_PS_Agree"
$ws.Range("A50").Value2 = "מאז הסקר האחרון, אבא שלי התעצבן או צעק עליי"
$ws.Range("B50").ClearContents()
$ws.Range("C50").Value2 = "C_PS_GotAngry"
$ws.Range("A51").Value2 = "מאז הסקר האחרון, אמא שלי התעצבנה או צעקה עליי"
$ws.Range("B51").ClearContents()
$ws.Range("C51").Value2 = "C_PS_GotAngry"
$ws.Range("A52").Value2 = "מאז הסקר האחרון, אבא שלי היה סבלני כלפיי"
$ws.Range("B52").ClearContents()
$ws.Range("C52").Value2 = "C_PS_Patient"
$ws.Range("A53").Value2 = "מאז הסקר האחרון, אמא שלי הייתה סבלנית כלפיי"
$ws.Range("B53").ClearContents()
$ws.Range("C53").Value2 = "C_PS_Patient"
$ws.Range("A54").Value2 = "תן דוגמה למשהו שעצבן אותך היום"
$ws.Range("B54").ClearContents()
$ws.Range("C54").Value2 = "C_triggers"
$ws.Range("A55").Value2 = "תני דוגמה למשהו שעצבן אותך היום"
$ws.Range("B55").ClearContents()
$ws.Range("C55").Value2 = "C_triggers"
# Bold the header row
$ws.Range("A1:C1").Font.Bold = $true

# Convert the data range into an Excel Table
$tableRange = $ws.Range("A1:C55")
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"

# Sort the table data by the label column (column C), ascending
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add2($ws.Range("C2:C55")) | Out-Null
$ws.Sort.SetRange($ws.Range("A1:C55"))
$ws.Sort.Header = 1
$ws.Sort.Apply()

# Update the active cell selection
$ws.Range("H13").Select()
